$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row right after row 9 (new row 10), shifting
#     everything below down by one. This is the "CAP2ACT is now entity
#     dependent" row added for conv_elec_hydroror (capacity_to_activity).
$ws.Rows("10:10").Insert()

$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "conv_elec_hydroror"
$ws.Range("C10").Value = "capacity_to_activity"
$ws.Range("D10").Value = "constant"
$ws.Range("G10").Value = 0.001
$ws.Range("H10").Value = "GW/TWh"

# --- Re-apply the AutoFilter so its range grows from L849 to L850 to
#     follow the extra row (the used range now ends one row later).
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A5:L850").AutoFilter()

# --- The hidden _FilterDatabase defined name tracks the same range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$850"
    }
}

# --- Move the active selection, matching the author's final cursor spot.
$ws.Range("B10").Select()

Write-Output "done"
